{"js": "// Insert a new \"List Bullet\" paragraph listing the responsible instructor\n// right after the \"Docente(s) Respons\u00e1vel(eis)\" heading paragraph.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst heading = paragraphs.items.find((p) =>\n  p.text.indexOf(\"Docente(s) Respons\u00e1vel(eis)\") !== -1\n);\n\nif (!heading) {\n  throw new Error('Could not find the \"Docente(s) Respons\u00e1vel(eis)\" heading paragraph.');\n}\n\nconst newParagraph = heading.insertParagraph(\n  \"6712818 - Mauricio Lamano Ferreira\",\n  \"After\"\n);\nnewParagraph.style = \"List Bullet\";\n\nawait context.sync();\n", "ps1": "# Insert a new \"List Bullet\" paragraph listing the responsible instructor\n# right after the \"Docente(s) Respons\u00e1vel(eis)\" heading paragraph.\n$d = $word.ActiveDocument\n\n$rng = $d.Content\n$found = $rng.Find.Execute(\"Docente(s) Respons\u00e1vel(eis)\")\nif (-not $found) {\n    throw 'Could not find the \"Docente(s) Respons\u00e1vel(eis)\" heading paragraph.'\n}\n\n$headingPara = $rng.Paragraphs(1)\n$insertionPoint = $headingPara.Range\n$insertionPoint.Collapse(0)  # wdCollapseEnd\n$insertionPoint.InsertParagraphAfter()\n\n$newPara = $headingPara.Next()\n$newPara.Range.Text = \"6712818 - Mauricio Lamano Ferreira\"\n$newPara.Style = \"List Bullet\"\n"}
